# Updated cryptos list on Wed Mar  6 12:12:17 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell in columns D (Price) and E (Volume(1h)) is stored as literal
# text in the workbook. A leading apostrophe forces Excel to keep the
# assigned value as text instead of coercing numeric-looking strings
# (e.g. "1.00", "0.0460") into Double/Int values, and resetting the
# range Style back to "Normal" afterwards drops the transient
# quote-prefix/text formatting so the cell keeps its original styling.

$ws.Range('D2').Value = "'" + '67.242.21'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +0.49%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '3.893.59'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  +3.65%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'" + '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.24%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '426.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  +1.71%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '131.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  -0.28%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'" + '3.890.52'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  +3.92%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'" + '0.614'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = "'" + '1.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  +0.07%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'" + '0.734'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -5.19%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  -7.70%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + '0.0000364'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  -10.00%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '40.89'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  -4.89%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '4.492.10'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  +3.88%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '10.07'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  -4.12%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '15.61'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  +18.47%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'" + '3.899.02'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  +4.78%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -1.20%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'" + '19.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  -5.52%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '67.555.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  +0.98%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -6.02%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '408.68'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  -8.39%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '14.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -12.25%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '85.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  -5.19%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '3.03'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  -4.24%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'" + '37.57'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  -3.43%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '5.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  +11.91%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'" + '3.21'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  -3.88%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -6.51%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'" + '694.94'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +4.81%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  -2.77%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'" + '12.44'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  -2.38%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  -0.35%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'" + '7.23'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  -1.02%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -8.66%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'" + '38.54'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  -8.75%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  +7.65%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -0.01%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'" + '55.44'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  -2.78%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '0.0460'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -6.43%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '3.01'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  -0.14%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  +0.42%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  -8.82%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'" + 'Monero'
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'" + 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'" + '147.49'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  +0.32%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = "'" + 'NEARProtocol'
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = "'" + 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = "'" + '4.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  +2.11%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '26.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  -8.02%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  -3.29%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'" + '3.27'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -5.27%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  -4.17%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -3.53%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '2.55'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  -4.58%  '
$ws.Range('E51').Style = 'Normal'
